$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B295").Value = 7047575
$ws.Range("F295").Value = "Crewe U21"
$ws.Range("G295").Value = "Sheff Utd U21"
$ws.Range("H295").Value = 2
$ws.Range("I295").Value = 5
$ws.Range("J295").Value = "A"
$ws.Range("K295").Value = 7
$ws.Range("L295").Value = 6
$ws.Range("M295").Value = 1.25
$ws.Range("N295").Value = 4
$ws.Range("O295").Value = 5
$ws.Range("P295").Value = 1.5
$ws.Range("Q295").Value = 1.25
$ws.Range("R295").Value = 1.85
$ws.Range("S295").Value = 1.95
$ws.Range("T295").Value = 3.75
$ws.Range("U295").Value = 1.9
$ws.Range("V295").Value = 1.9
$ws.Range("W295").Value = -1
$ws.Range("X295").Value = -1
$ws.Range("Y295").Value = 0.5
$ws.Range("Z295").Value = -1
$ws.Range("AA295").Value = 0.95
$ws.Range("AB295").Value = 0.8999999999999999
$ws.Range("AC295").Value = -1

$ws.Range("B296").Value = 7047576
$ws.Range("F296").Value = "Wigan U21"
$ws.Range("G296").Value = "Birmingham U21"
$ws.Range("H296").Value = 0
$ws.Range("I296").Value = 5
$ws.Range("J296").Value = "A"
$ws.Range("K296").Value = 2
$ws.Range("L296").Value = 4
$ws.Range("M296").Value = 2.75
$ws.Range("N296").Value = 2.875
$ws.Range("O296").Value = 4
$ws.Range("P296").Value = 1.95
$ws.Range("Q296").Value = 0.5
$ws.Range("R296").Value = 1.8
$ws.Range("S296").Value = 2
$ws.Range("T296").Value = 3.5
$ws.Range("U296").Value = 1.925
$ws.Range("V296").Value = 1.875
$ws.Range("W296").Value = -1
$ws.Range("X296").Value = -1
$ws.Range("Y296").Value = 0.95
$ws.Range("Z296").Value = -1
$ws.Range("AA296").Value = 1
$ws.Range("AB296").Value = 0.925
$ws.Range("AC296").Value = -1

$ws.Range("B313").Value = 7047602
$ws.Range("F313").Value = "Coventry U21"
$ws.Range("G313").Value = "Crewe U21"
$ws.Range("H313").Value = 2
$ws.Range("I313").Value = 0
$ws.Range("J313").Value = "H"
$ws.Range("K313").Value = 1.5
$ws.Range("L313").Value = 4.5
$ws.Range("M313").Value = 4.5
$ws.Range("N313").Value = 1.75
$ws.Range("O313").Value = 4.2
$ws.Range("P313").Value = 3.4
$ws.Range("Q313").Value = -0.75
$ws.Range("R313").Value = 1.975
$ws.Range("S313").Value = 1.825
$ws.Range("T313").Value = 3.5
$ws.Range("U313").Value = 1.775
$ws.Range("V313").Value = 1.925
$ws.Range("W313").Value = 0.75
$ws.Range("X313").Value = -1
$ws.Range("Y313").Value = -1
$ws.Range("Z313").Value = 0.9750000000000001
$ws.Range("AA313").Value = -1
$ws.Range("AB313").Value = -1
$ws.Range("AC313").Value = 0.925

$ws.Range("B314").Value = 7047710
$ws.Range("F314").Value = "Ipswich U21"
$ws.Range("G314").Value = "Charlton U21"
$ws.Range("H314").Value = 1
$ws.Range("I314").Value = 8
$ws.Range("J314").Value = "A"
$ws.Range("K314").Value = 2.25
$ws.Range("L314").Value = 4
$ws.Range("M314").Value = 2.4
$ws.Range("N314").Value = 2.6
$ws.Range("O314").Value = 3.8
$ws.Range("P314").Value = 2.15
$ws.Range("Q314").Value = 0.25
$ws.Range("R314").Value = 1.825
$ws.Range("S314").Value = 1.975
$ws.Range("T314").Value = 3.5
$ws.Range("U314").Value = 1.825
$ws.Range("V314").Value = 1.975
$ws.Range("W314").Value = -1
$ws.Range("X314").Value = -1
$ws.Range("Y314").Value = 1.15
$ws.Range("Z314").Value = -1
$ws.Range("AA314").Value = 0.9750000000000001
$ws.Range("AB314").Value = 0.825
$ws.Range("AC314").Value = -1

$ws.Range("B315").Value = 7047603
$ws.Range("F315").Value = "Sheffield Wed U21"
$ws.Range("G315").Value = "Burnley U21"
$ws.Range("H315").Value = 2
$ws.Range("I315").Value = 1
$ws.Range("J315").Value = "H"
$ws.Range("K315").Value = 1.9
$ws.Range("L315").Value = 4
$ws.Range("M315").Value = 3
$ws.Range("N315").Value = 1.7
$ws.Range("O315").Value = 4.2
$ws.Range("P315").Value = 3.6
$ws.Range("Q315").Value = -0.75
$ws.Range("R315").Value = 1.925
$ws.Range("S315").Value = 1.875
$ws.Range("T315").Value = 3.5
$ws.Range("U315").Value = 1.975
$ws.Range("V315").Value = 1.825
$ws.Range("W315").Value = 0.7
$ws.Range("X315").Value = -1
$ws.Range("Y315").Value = -1
$ws.Range("Z315").Value = 0.4625
$ws.Range("AA315").Value = -0.5
$ws.Range("AB315").Value = -1
$ws.Range("AC315").Value = 0.825
